# Weekly update: a new "Ajo" (garlic) price entry for Femacal de La Calera
# was reported, pushing all the existing rows from 208 downward by one row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row at 208; Excel shifts rows 208:244 down to 209:245.
$ws.Rows("208:208").Insert()

# Populate the newly inserted row with the new reported price observation.
$ws.Range("A208").Value = 3
$ws.Range("B208").Value = "Femacal de La Calera"
$ws.Range("C208").Value = "Coquimbo"
$ws.Range("D208").Value = 44504
$ws.Range("E208").Value = 5
$ws.Range("F208").Value = 100112003
$ws.Range("G208").Value = "Ajo"
$ws.Range("H208").Value = "Chino"
$ws.Range("I208").Value = "Primera"
$ws.Range("J208").Value = 73
$ws.Range("K208").Value = 16000
$ws.Range("L208").Value = 16500
$ws.Range("M208").Value = 16260
$ws.Range("N208").Value = "$/caja 10 kilos"
$ws.Range("O208").Value = "China"
$ws.Range("P208").Value = 1626
$ws.Range("Q208").Value = 10
$ws.Range("R208").Value = "Hortaliza"
